# Update the title on the "SERVEUR - API - Arduino" slide so that the
# trailing word becomes "ARDUINO" and the text ends up split across three
# runs (matching how PowerPoint splits runs when text is edited in place):
#   "SERVEUR – API "  +  "– "  +  "ARDUINO"

$p = $ppt.ActivePresentation

$dash = [string][char]0x2013   # U+2013 EN DASH ("–") used in the original text

# Locate the slide/shape that holds the title text, rather than hard-coding
# indices, so the script is resilient to ordering quirks.
$targetSlide = $null
$targetShape = $null

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame) {
            if ($shape.TextFrame.TextRange.Text -like "*Arduino*") {
                $targetSlide = $slide
                $targetShape = $shape
            }
        }
    }
}

$tr = $targetShape.TextFrame.TextRange

# Step 1: split off the "- " (dash + space) right before "Arduino" into its
# own run by re-typing it.
$fullText = $tr.Text
$dashIdx = $fullText.IndexOf($dash + " Arduino")
$dashRange = $tr.Characters($dashIdx + 1, 2)
$dashRange.Text = $dash + " "

# Step 2: upper-case "Arduino" -> "ARDUINO", which also splits it into its
# own trailing run.
$fullText2 = $tr.Text
$arduinoIdx = $fullText2.IndexOf("Arduino")
$arduinoRange = $tr.Characters($arduinoIdx + 1, 7)
$arduinoRange.Text = "ARDUINO"

Write-Host "Updated title text: $($tr.Text)"
Write-Host "Run count: $($tr.Runs().Count)"
